$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1550.186
$ws.Range("I15").Value = 1550.186
$ws.Range("K15").Value = 4650.558
$ws.Range("M15").Value = -4481.558
$ws.Range("H47").Value = 7000
$ws.Range("I47").Value = 6000
$ws.Range("K47").Value = 6000
$ws.Range("M47").Value = -5028
$ws.Range("H86").Value = 7713.7
$ws.Range("I86").Value = 2639.9167
$ws.Range("J86").Value = 15324.375
$ws.Range("K86").Value = 2639.9167
$ws.Range("L86").Value = 15324.375
$ws.Range("M86").Value = -1516.9167
$ws.Range("N86").Value = -17570.375
$ws.Range("H89").Value = 7713.7
$ws.Range("I89").Value = 2639.9167
$ws.Range("J89").Value = 15324.375
$ws.Range("K89").Value = 13199.5835
$ws.Range("L89").Value = 76621.875
$ws.Range("M89").Value = -7583.583500000001
$ws.Range("N89").Value = -87853.875
$ws.Range("H112").Value = 1086.742
$ws.Range("J112").Value = 1086.742
$ws.Range("L112").Value = 3260.226
$ws.Range("N112").Value = -5476.226
$ws.Range("H129").Value = 2458.739
$ws.Range("J129").Value = 3233.5293
$ws.Range("L129").Value = 9700.5879
$ws.Range("N129").Value = -19700.5879
$ws.Range("H137").Value = 2581.0908
$ws.Range("J137").Value = 2699
$ws.Range("L137").Value = 8097
$ws.Range("N137").Value = -13197
$ws.Range("H138").Value = 2099.0125
$ws.Range("J138").Value = 2205.4355
$ws.Range("L138").Value = 6616.306500000001
$ws.Range("N138").Value = -16896.3065
$ws.Range("H141").Value = 1102.8064
$ws.Range("I141").Value = 1006.23334
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 3018.70002
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 2161.29998
$ws.Range("N141").Value = -22360

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4977.356
$ws.Range("I32").Value = 5132.46
$ws.Range("J32").Value = 4115.6665
$ws.Range("K32").Value = 5132.46
$ws.Range("L32").Value = 4115.6665
$ws.Range("M32").Value = -4845.46
$ws.Range("N32").Value = -4689.6665
$ws.Range("H74").Value = 125000800
$ws.Range("I74").Value = 166667460
$ws.Range("J74").Value = 807
$ws.Range("K74").Value = 166667460
$ws.Range("L74").Value = 807
$ws.Range("M74").Value = -166666586
$ws.Range("N74").Value = -2555
$ws.Range("H77").Value = 125000800
$ws.Range("I77").Value = 166667460
$ws.Range("J77").Value = 807
$ws.Range("K77").Value = 833337300
$ws.Range("L77").Value = 4035
$ws.Range("M77").Value = -833332932
$ws.Range("N77").Value = -12771
$ws.Range("H88").Value = 500527
$ws.Range("I88").Value = 1000
$ws.Range("K88").Value = 1000
$ws.Range("M88").Value = -594
$ws.Range("H91").Value = 500527
$ws.Range("I91").Value = 1000
$ws.Range("K91").Value = 1000
$ws.Range("M91").Value = 404
$ws.Range("H110").Value = 534.1
$ws.Range("J110").Value = 546.4
$ws.Range("L110").Value = 546.4
$ws.Range("N110").Value = -4636.4

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1580.3077
$ws.Range("I20").Value = 1873.8889
$ws.Range("J20").Value = 919.75
$ws.Range("K20").Value = 1873.8889
$ws.Range("L20").Value = 919.75
$ws.Range("M20").Value = -1626.8889
$ws.Range("N20").Value = -1413.75
$ws.Range("H105").Value = 4042.4546
$ws.Range("I105").Value = 4744.6665
$ws.Range("K105").Value = 4744.6665
$ws.Range("M105").Value = -2997.6665

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15828.6
$ws.Range("I31").Value = 24752.53
$ws.Range("J31").Value = 4158.846
$ws.Range("K31").Value = 24752.53
$ws.Range("L31").Value = 4158.846
$ws.Range("M31").Value = -24457.53
$ws.Range("N31").Value = -4748.846
$ws.Range("H34").Value = 15828.6
$ws.Range("I34").Value = 24752.53
$ws.Range("J34").Value = 4158.846
$ws.Range("K34").Value = 24752.53
$ws.Range("L34").Value = 4158.846
$ws.Range("M34").Value = -24550.53
$ws.Range("N34").Value = -4562.846
$ws.Range("H60").Value = 14650
$ws.Range("I60").Value = 3000
$ws.Range("J60").Value = 18533.334
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 18533.334
$ws.Range("M60").Value = -2489
$ws.Range("N60").Value = -19555.334
$ws.Range("H74").Value = 30296.727
$ws.Range("J74").Value = 30296.727
$ws.Range("L74").Value = 30296.727
$ws.Range("N74").Value = -32044.727
$ws.Range("H77").Value = 30296.727
$ws.Range("J77").Value = 30296.727
$ws.Range("L77").Value = 90890.181
$ws.Range("N77").Value = -99626.181

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 3332
$ws.Range("J58").Value = 3720.6667
$ws.Range("L58").Value = 11162.0001
$ws.Range("N58").Value = -11418.0001
$ws.Range("H68").Value = 17250.5
$ws.Range("I68").Value = 750
$ws.Range("J68").Value = 25500.75
$ws.Range("K68").Value = 2250
$ws.Range("L68").Value = 76502.25
$ws.Range("M68").Value = -1439
$ws.Range("N68").Value = -78124.25
$ws.Range("H69").Value = 2500
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2500
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 7500
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -9122
$ws.Range("H71").Value = 17250.5
$ws.Range("I71").Value = 750
$ws.Range("J71").Value = 25500.75
$ws.Range("K71").Value = 6750
$ws.Range("L71").Value = 229506.75
$ws.Range("M71").Value = -2694
$ws.Range("N71").Value = -237618.75
$ws.Range("H72").Value = 2500
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2500
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 22500
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -30612
$ws.Range("H113").Value = 906.6667
$ws.Range("I113").Value = 900
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = -530
$ws.Range("H122").Value = 693.7083
$ws.Range("J122").Value = 891.8
$ws.Range("L122").Value = 8026.2
$ws.Range("N122").Value = -12926.2
$ws.Range("H131").Value = 811.14
$ws.Range("J131").Value = 811.2525
$ws.Range("L131").Value = 2433.7575
$ws.Range("N131").Value = -12513.7575

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 332.33334
$ws.Range("I9").Value = 332.33334
$ws.Range("K9").Value = 332.33334
$ws.Range("M9").Value = -162.33334
$ws.Range("H80").Value = 3839.1875
$ws.Range("J80").Value = 3918.5
$ws.Range("L80").Value = 3918.5
$ws.Range("N80").Value = -5914.5
$ws.Range("H83").Value = 3839.1875
$ws.Range("J83").Value = 3918.5
$ws.Range("L83").Value = 19592.5
$ws.Range("N83").Value = -29576.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3018.9167
$ws.Range("I40").Value = 2177.75
$ws.Range("K40").Value = 2177.75
$ws.Range("M40").Value = -2041.75
$ws.Range("H82").Value = 1875.0333
$ws.Range("I82").Value = 1762.04
$ws.Range("J82").Value = 2440
$ws.Range("K82").Value = 1762.04
$ws.Range("L82").Value = 2440
$ws.Range("M82").Value = -1401.04
$ws.Range("N82").Value = -3162
$ws.Range("H85").Value = 1875.0333
$ws.Range("I85").Value = 1762.04
$ws.Range("J85").Value = 2440
$ws.Range("K85").Value = 1762.04
$ws.Range("L85").Value = 2440
$ws.Range("M85").Value = -514.04
$ws.Range("N85").Value = -4936
$ws.Range("H122").Value = 1034828.1
$ws.Range("I122").Value = 2453804.2
$ws.Range("K122").Value = 7361412.600000001
$ws.Range("M122").Value = -7358962.600000001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3333.1667
$ws.Range("J15").Value = 3333.1667
$ws.Range("L15").Value = 3333.1667
$ws.Range("N15").Value = -3909.1667
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H69").Value = 16911.5
$ws.Range("J69").Value = 15800
$ws.Range("L69").Value = 15800
$ws.Range("N69").Value = -17298
$ws.Range("H72").Value = 16911.5
$ws.Range("J72").Value = 15800
$ws.Range("L72").Value = 47400
$ws.Range("N72").Value = -54888
$ws.Range("H81").Value = 100001870
$ws.Range("I81").Value = 1780.375
$ws.Range("J81").Value = 500002240
$ws.Range("K81").Value = 3560.75
$ws.Range("L81").Value = 1000004480
$ws.Range("M81").Value = -2499.75
$ws.Range("N81").Value = -1000006602
$ws.Range("H84").Value = 100001870
$ws.Range("I84").Value = 1780.375
$ws.Range("J84").Value = 500002240
$ws.Range("K84").Value = 17803.75
$ws.Range("L84").Value = 5000022400
$ws.Range("M84").Value = -12499.75
$ws.Range("N84").Value = -5000033008
$ws.Range("H100").Value = 472.875
$ws.Range("J100").Value = 500
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082
